$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump "Last Updated" timestamp by a minute ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 06:18 PM"

# --- 2. Top Gainers sheet: refresh figures, drop GREENLAM, add CGPOWER, ---
#        then re-sort the leaderboard by "Latest" (column C), descending,
#        same as the live sheet does whenever the feed updates.
$wsGainers = $wb.Worksheets.Item("Top Gainers")

# KERNEX (row 57) gets refreshed Latest/Weekly/Monthly figures.
$wsGainers.Range("C57").Value = 3.9981
$wsGainers.Range("D57").Value = 7.4592
$wsGainers.Range("E57").Value = 27.1054

# NPST (row 67) gets refreshed Latest/Weekly/Monthly figures too - this
# moves it up above ORIENTTECH/ICRA/SALASAR once the list is re-sorted.
$wsGainers.Range("C67").Value = 3.8509
$wsGainers.Range("D67").Value = -2.0059
$wsGainers.Range("E67").Value = -3.5057

# GREENLAM (row 61) drops off the list entirely; everything below shifts up.
$wsGainers.Rows.Item(61).Delete()

# A new entry, CGPOWER, joins the bottom of the table (row 76).
$wsGainers.Range("A76").Value = "🚀"
$wsGainers.Range("B76").Value = "CGPOWER"
$wsGainers.Range("C76").Value = 3.6125
$wsGainers.Range("D76").Value = 3.4192
$wsGainers.Range("E76").Value = 1.0325

# Re-sort the data rows (A2:E76) by column C, descending, so the
# leaderboard order reflects the refreshed "Latest" figures.
$sortRange = $wsGainers.Range("A2:E76")
$sortRange.Sort($wsGainers.Range("C2:C76"), 2, $null, $null, 1, $null, 1, 0)

# --- 3. Top Losers sheet: update a few Weekly (column D) figures ---
$wsLosers = $wb.Worksheets.Item("Top Losers")

# CRAMC (row 18)
$wsLosers.Range("D18").Value = -0.062

# RUBICON (row 48)
$wsLosers.Range("D48").Value = 0.05

# CANHLIFE (row 54) now has a real Weekly figure instead of "N/A"
$wsLosers.Range("D54").Value = 5.2953
